$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet: "Voice Lines - main" -> "Text Lines - main"
$ws.Name = "Text Lines - main"

# Remove the rows that correspond to deleted scene/test lines.
# (Original row numbers, deleted from bottom to top so indices stay valid.)
$rowsToDelete = 29, 24, 23, 16, 14, 5, 4, 3
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
